$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.710.44'
$ws.Range("E2").Value = '  +4.04%  '

$ws.Range("D3").Value = '1.924.05'
$ws.Range("E3").Value = '  +2.76%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -1.21%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '335.73'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.10%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.46%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4694'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.36%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4139'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.63%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.20'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.92%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08053'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.79%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.017'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.02%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.41'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.76%  '

$ws.Range("D13").Value = '1.942.64'
$ws.Range("E13").Value = '  +2.78%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.028'
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.207'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.33%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '90.10'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.80%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.001'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.19%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001035'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.15%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06592'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.25%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.89'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.00%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9989'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.17%  '

$ws.Range("D22").Value = '29.667.30'
$ws.Range("E22").Value = '  +3.98%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.566'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.13%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.63'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +7.93%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.203'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.15%  '

$ws.Range("D26").Value = '2.144.41'
$ws.Range("E26").Value = '  +0.38%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '156.75'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.63%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.93'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.55%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.149'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.88%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.748'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +7.46%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '117.85'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.10%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.060'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +11.49%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09475'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.94%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.441'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.02%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.438'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.15%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.524'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.86%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06159'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.79%  '

$ws.Range("E38").Value = '  +2.50%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.465'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.03%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.187'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.44%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5922'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.28%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1849'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.66%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '10.29'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.62%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.258'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.68%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.348'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.62%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.07526'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.71%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5601'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.40%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '12.21'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.13%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.945'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.71%  '

$ws.Range("E50").Value = '  +2.26%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3001'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +12.48%  '
